# Add two new columns (I: "I0", J: "IF") to the sheet, with header cells
# in row 1 styled like the existing header cells, and fill in the values
# for rows 2-35 exactly as described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing last header cell (H1) onto the two
# new header cells so they share the same bold/border/center style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Header row (row 1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-35
$data = @(
    @(7, 8),
    @(8, 8),
    @(8, 9),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(4, 5),
    @(8, 8),
    @(6, 7),
    @(7, 8),
    @(7, 7),
    @(7, 7),
    @(8, 9),
    @(7, 7),
    @(10, 10),
    @(7, 8),
    @(7, 7),
    @(8, 8),
    @(6, 6),
    @(8, 8),
    @(6, 6),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(6, 7),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(4, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
